{"js": "// Acceptance Test Report.docx edit\n// 1) Title: \"Acceptance Testing\" -> \"Acceptance Test\"\n// 2) Move the \"_GoBack\" bookmark from the \"Not implemented.\" paragraph\n//    (end of \"Pick up assigned manuscripts\") to the new sentence appended\n//    at the very end of the document.\n// 3) Append \"The tests were run using Chrome, Firefox and Safari.\" to the\n//    last (previously empty) paragraph, removing its list indentation.\n\nconst body = context.document.body;\n\n// --- 1. Fix the document title -------------------------------------------------\nconst titleResults = body.search(\"Acceptance Testing\", { matchCase: true });\ntitleResults.load(\"text\");\nawait context.sync();\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\"Acceptance Test\", \"Replace\");\n}\nawait context.sync();\n\n// --- 2. Drop the old _GoBack bookmark ------------------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 3. Fill in the final paragraph --------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.leftIndent = 0;\nlastParagraph.insertText(\n  \"The tests were run using Chrome, Firefox and Safari.\",\n  \"Start\"\n);\nawait context.sync();\n\n// Re-create the _GoBack bookmark right after \"The tests were \" (collapsed,\n// i.e. zero-length) \u2014 this is where Word leaves it after the cursor last\n// stopped there while typing the sentence.\nconst wereResults = body.search(\"were \", { matchCase: true });\nwereResults.load(\"text\");\nawait context.sync();\nconst afterWere = wereResults.items[0].getRange(\"After\");\nafterWere.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// Force the same run split Word produces between \"run\" and \" using Chrome,\n// Firefox and Safari.\" (separate typing sessions become separate <w:r>\n// elements) by dropping a transient bookmark at that boundary and removing\n// it again immediately \u2014 inserting/removing a bookmark leaves the run break\n// behind without leaving any bookmark residue.\nconst runResults = body.search(\"were run\", { matchCase: true });\nrunResults.load(\"text\");\nawait context.sync();\nconst afterRun = runResults.items[0].getRange(\"After\");\nafterRun.insertBookmark(\"__tmp_run_split__\");\nawait context.sync();\n\ncontext.document.deleteBookmark(\"__tmp_run_split__\");\nawait context.sync();\n", "ps1": "# Acceptance Test Report.docx edit\n# 1) Title: \"Acceptance Testing\" -> \"Acceptance Test\"\n# 2) Move the \"_GoBack\" bookmark from the \"Not implemented.\" paragraph\n#    (end of \"Pick up assigned manuscripts\") to the new sentence appended\n#    at the very end of the document.\n# 3) Append \"The tests were run using Chrome, Firefox and Safari.\" to the\n#    last (previously empty) paragraph, removing its list indentation.\n\n$d = $word.ActiveDocument\n\n# --- 1. Fix the document title -------------------------------------------------\n$d.Content.Find.Execute(\"Acceptance Testing\", $false, $false, $false, $false, $false, $true, 1, $false, \"Acceptance Test\", 2) | Out-Null\n\n# --- 2. Drop the old _GoBack bookmark ------------------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- 3. Fill in the final paragraph --------------------------------------------\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertBefore(\"The tests were run using Chrome, Firefox and Safari.\")\n$lastParagraph.Format.LeftIndent = 0\n\n# Re-create the _GoBack bookmark right after \"The tests were \" (collapsed,\n# i.e. zero-length) -- this is where Word leaves it after the cursor last\n# stopped there while typing the sentence.\n$wereRange = $d.Content\n$wereRange.Find.Execute(\"were \")\n$bookmarkPoint = $wereRange.Duplicate\n$bookmarkPoint.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkPoint)\n\n# Force the same run split Word produces between \"run\" and \" using Chrome,\n# Firefox and Safari.\" (separate typing sessions become separate runs) by\n# dropping a transient bookmark at that boundary and removing it again\n# immediately -- inserting/removing a bookmark leaves the run break behind\n# without leaving any bookmark residue.\n$runRange = $d.Content\n$runRange.Find.Execute(\"were run\")\n$splitPoint = $runRange.Duplicate\n$splitPoint.Collapse(0)\n$d.Bookmarks.Add(\"__tmp_run_split__\", $splitPoint)\n$d.Bookmarks(\"__tmp_run_split__\").Delete()\n"}
